# Introduction to Stratoshark - trim the demo lineup down to a single demo.
#
# 1. Rename "Demo 1: HTTP" -> "Demo: HTTP" (slide with that title, currently
#    slide 11 of the deck).
# 2. Delete the "Demo 2: SCP File Transfer" slide entirely (currently slide
#    12 of the deck, immediately after the HTTP demo slide).

$p = $ppt.ActivePresentation

# --- Step 1: rename the remaining demo slide's title -----------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq "Demo 1: HTTP") {
                $shape.TextFrame.TextRange.Text = "Demo: HTTP"
            }
        }
    }
}

# --- Step 2: delete the "Demo 2: SCP File Transfer" slide -------------------
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $isDemo2 = $false
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq "Demo 2: SCP File Transfer") {
                $isDemo2 = $true
            }
        }
    }
    if ($isDemo2) {
        $slide.Delete()
    }
}
